# Apply updated coin price / volume(1h) values scraped on
# Thu Jan 26 22:39:39 UTC 2023 by the symbol-list GitHub Action.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each cell originally holds a plain text (inline string) value - e.g. "304.96"
# or "-0.92%". Assigning a bare numeric-looking string via .Value lets Excel
# auto-convert it to a real number (and percentages to fractional numbers), so
# we force text with a leading apostrophe and then reset the style back to
# "Normal" so no stray Text-format style is left behind on the cell.

function Set-TextValue($cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" "305.19"
Set-TextValue "E2" "-0.64%"
Set-TextValue "E3" "-1.93%"
Set-TextValue "D4" "4.985"
Set-TextValue "E4" "-1.73%"
Set-TextValue "E5" "-1.54%"
Set-TextValue "D6" "1.886"
Set-TextValue "E6" "-4.05%"
Set-TextValue "D7" "4.148"
Set-TextValue "E7" "1.54%"
Set-TextValue "D8" "7.868"
Set-TextValue "E8" "0.50%"
Set-TextValue "D9" "0.9297"
Set-TextValue "E9" "-0.54%"
Set-TextValue "D10" "0.1326"
Set-TextValue "E10" "-10.19%"
Set-TextValue "D11" "0.1904"
Set-TextValue "E11" "-1.64%"
Set-TextValue "D12" "0.09233"
Set-TextValue "E12" "1.19%"
Set-TextValue "D13" "0.03511"
Set-TextValue "E13" "-0.22%"
Set-TextValue "D14" "0.09914"
Set-TextValue "E14" "1.05%"
Set-TextValue "D15" "0.001414"
Set-TextValue "E15" "0.07%"
Set-TextValue "D16" "0.006351"
Set-TextValue "E16" "10.03%"
Set-TextValue "D17" "3.606"
Set-TextValue "E17" "2.31%"
Set-TextValue "E19" "0.82%"
Set-TextValue "D20" "5.237"
Set-TextValue "E20" "3.72%"
Set-TextValue "E21" "0.07%"
Set-TextValue "E22" "5.74%"
Set-TextValue "D23" "0.04405"
Set-TextValue "E23" "-2.07%"
Set-TextValue "E24" "2.17%"
Set-TextValue "D25" "0.004728"
Set-TextValue "E25" "-3.64%"
Set-TextValue "E26" "5.97%"
Set-TextValue "D27" "0.0003129"
Set-TextValue "E27" "-29.36%"
Set-TextValue "D39" "0.01956"
Set-TextValue "E39" "-1.67%"
Set-TextValue "D40" "0.05225"
Set-TextValue "E40" "7.54%"
Set-TextValue "D41" "0.007548"
Set-TextValue "E41" "0.09%"
Set-TextValue "D42" "0.01017"
Set-TextValue "E42" "-7.63%"
Set-TextValue "D43" "0.1373"
Set-TextValue "E43" "-0.57%"
Set-TextValue "E44" "1.38%"
Set-TextValue "D45" "0.01072"
Set-TextValue "E45" "-1.71%"
Set-TextValue "D46" "0.00006330"
Set-TextValue "E46" "3.32%"
Set-TextValue "E47" "0.29%"
Set-TextValue "D48" "64.96"
Set-TextValue "E48" "0.45%"
Set-TextValue "D49" "0.001659"
Set-TextValue "E49" "39.74%"
Set-TextValue "E50" "0.29%"
Set-TextValue "E51" "0.29%"
